# Apply the 2024-11-28 cryptos list refresh (prices + 1h volume % changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) sometimes hold plain decimal-looking text (e.g. "237.64").
# Force those specific cells to Text format first so Excel keeps them as
# strings instead of silently coercing them into floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Write the updated Price (D) and Volume(1h) (E) values for each row.
$ws.Range("D2").Value = "95.306.62"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "3.586.46"
$ws.Range("E3").Value = "  +4.74%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "237.64"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "649.40"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("D7").Value = "1.46"
$ws.Range("E7").Value = "  +3.54%  "
$ws.Range("D8").Value = "0.401"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "0.992"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").Value = "3.585.17"
$ws.Range("E11").Value = "  +4.69%  "
$ws.Range("D12").Value = "42.55"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "0.199"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "6.30"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "4.269.11"
$ws.Range("E15").Value = "  +4.98%  "
$ws.Range("D16").Value = "95.311.05"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "0.0000253"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "3.589.55"
$ws.Range("E18").Value = "  +4.83%  "
$ws.Range("D19").Value = "7.91"
$ws.Range("E19").Value = "  -4.49%  "
$ws.Range("D20").Value = "12.66"
$ws.Range("E20").Value = "  +8.42%  "
$ws.Range("D21").Value = "17.86"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").Value = "3.61"
$ws.Range("E22").Value = "  +6.35%  "
$ws.Range("D23").Value = "0.485"
$ws.Range("E23").Value = "  +3.10%  "
$ws.Range("D24").Value = "507.70"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").Value = "0.0000194"
$ws.Range("E25").Value = "  +4.56%  "
$ws.Range("D26").Value = "6.56"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").Value = "96.23"
$ws.Range("E27").Value = "  +4.19%  "
$ws.Range("D28").Value = "12.60"
$ws.Range("E28").Value = "  +4.66%  "
$ws.Range("D29").Value = "3.779.76"
$ws.Range("E29").Value = "  +4.65%  "
$ws.Range("D30").Value = "3.08"
$ws.Range("E30").Value = "  +12.00%  "
$ws.Range("D31").Value = "11.28"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "0.139"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "0.986"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "0.177"
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("D36").Value = "31.67"
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("D37").Value = "0.556"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "8.13"
$ws.Range("E38").Value = "  +8.21%  "
$ws.Range("D39").Value = "569.57"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "0.915"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").Value = "1.71"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "5.65"
$ws.Range("E46").Value = "  +2.19%  "
$ws.Range("D47").Value = "33.98"
$ws.Range("E47").Value = "  +35.78%  "
$ws.Range("D48").Value = "2.23"
$ws.Range("E48").Value = "  +4.36%  "
$ws.Range("D49").Value = "0.0412"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").Value = "3.52"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("D51").Value = "53.72"
$ws.Range("E51").Value = "  +0.19%  "
